$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These titles were dropped from the PBC study; their address column is
# marked as "not included" (plain text, no hyperlink / hyperlink style)
# instead of the PBC URL that used to live there.
$targets = @('$B$15', '$B$58', '$B$91', '$B$104')

foreach ($h in @($ws.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($targets -contains $addr) {
        $h.Delete()
    }
}

$rows = @(15, 58, 91, 104)
foreach ($r in $rows) {
    $cell = $ws.Range("B$r")
    $cell.ClearFormats()
    $cell.Value = "NIEUWZGLĘDNIONE"
}

# Move the active selection / viewport.
$ws.Range("E110").Select()
